$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (Changed) date column C for all existing data
#    rows (2 through 498) from 2023-09-06 (45175) to 2023-09-08 (45177).
$ws.Range("C2:C498").Value2 = 45177

# 2) The existing last data row (498) now gets an explicit row height.
$ws.Rows.Item(498).RowHeight = 15

# 3) Append new row 499: "A 41425-2023"
$ws.Rows.Item(499).RowHeight = 15

$ws.Cells.Item(499, 1).Value = "A 41425-2023"

$ws.Cells.Item(499, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(499, 2).Value2 = 45175

$ws.Cells.Item(499, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(499, 3).Value2 = 45177

$ws.Cells.Item(499, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(499, 5).Value = "FALKENBERG"

$ws.Cells.Item(499, 7).Value2 = 0.6
$ws.Cells.Item(499, 8).Value2 = 0
$ws.Cells.Item(499, 9).Value2 = 0
$ws.Cells.Item(499, 10).Value2 = 0
$ws.Cells.Item(499, 11).Value2 = 0
$ws.Cells.Item(499, 12).Value2 = 0
$ws.Cells.Item(499, 13).Value2 = 0
$ws.Cells.Item(499, 14).Value2 = 0
$ws.Cells.Item(499, 15).Value2 = 0
$ws.Cells.Item(499, 16).Value2 = 0
$ws.Cells.Item(499, 17).Value2 = 0

$ws.Cells.Item(499, 18).WrapText = $true

# 4) Append new row 500: "A 41835-2023"
$ws.Cells.Item(500, 1).Value = "A 41835-2023"

$ws.Cells.Item(500, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(500, 2).Value2 = 45176

$ws.Cells.Item(500, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(500, 3).Value2 = 45177

$ws.Cells.Item(500, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(500, 5).Value = "FALKENBERG"

$ws.Cells.Item(500, 7).Value2 = 8.9
$ws.Cells.Item(500, 8).Value2 = 0
$ws.Cells.Item(500, 9).Value2 = 0
$ws.Cells.Item(500, 10).Value2 = 0
$ws.Cells.Item(500, 11).Value2 = 0
$ws.Cells.Item(500, 12).Value2 = 0
$ws.Cells.Item(500, 13).Value2 = 0
$ws.Cells.Item(500, 14).Value2 = 0
$ws.Cells.Item(500, 15).Value2 = 0
$ws.Cells.Item(500, 16).Value2 = 0
$ws.Cells.Item(500, 17).Value2 = 0

$ws.Cells.Item(500, 18).WrapText = $true
